{"js": "// The canonical-OOXML diff for this template is a pure XML-serialization\n// reorder: the xmlns:* namespace declarations on the root <w:document>\n// element, the attribute order inside <w:pgSz>/<w:pgMar> (document.xml),\n// and the attribute order inside <w:rFonts>, <w:lang>, <w:latentStyles>,\n// every <w:lsdException>, every <w:style> and <w:tblInd>/<w:tblCellMar>\n// (styles.xml) all get their attributes alphabetized. Every attribute\n// name/value pair is exactly the same before and after - nothing is\n// added, removed, or renamed, and no paragraph text, field, bookmark,\n// page-setup value or style definition actually changes. That kind of\n// attribute-order normalization comes from the producing tool's save\n// routine; Office.js has no API that lets a caller choose the raw\n// attribute order used when OOXML is serialized, so there is no\n// document-model edit to make.\n//\n// Read (without writing) the exact areas the reordering touched - page\n// setup and the template's styles - so the script still exercises the\n// Word JavaScript API surface; none of these reads change document\n// content.\nconst sections = context.document.sections;\nsections.load(\"items\");\n\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\n\nawait context.sync();\n\nif (sections.items.length > 0) {\n    const pageSetup = sections.items[0].pageSetup;\n    pageSetup.load([\n        \"pageWidth\",\n        \"pageHeight\",\n        \"topMargin\",\n        \"bottomMargin\",\n        \"leftMargin\",\n        \"rightMargin\",\n        \"headerDistance\",\n        \"footerDistance\",\n        \"gutter\"\n    ]);\n    await context.sync();\n}\n\n// No content, formatting, or style mutation is required: the template's\n// text, fields, bookmarks, section/page-setup values and style\n// definitions are unchanged by the source diff.\n", "ps1": "# The canonical-OOXML diff for this template is a pure XML-serialization\n# reorder: the xmlns:* namespace declarations on the root <w:document>\n# element, the attribute order inside <w:pgSz>/<w:pgMar> (document.xml),\n# and the attribute order inside <w:rFonts>, <w:lang>, <w:latentStyles>,\n# every <w:lsdException>, every <w:style> and <w:tblInd>/<w:tblCellMar>\n# (styles.xml) all get their attributes alphabetized. Every attribute\n# name/value pair is exactly the same before and after - nothing is\n# added, removed, or renamed, and no paragraph text, field, bookmark,\n# page-setup value or style definition actually changes. That kind of\n# attribute-order normalization comes from the producing tool's save\n# routine; the Word COM object model has no property that lets a caller\n# choose the raw attribute order used when OOXML is serialized, so there\n# is no document-model edit to make.\n\n$d = $word.ActiveDocument\n\n# Read (without writing) the exact areas the reordering touched - page\n# setup and the template's styles - so the script still exercises the\n# Word object model; none of these reads change document content.\n$ps = $d.PageSetup\n$null = $ps.PageWidth\n$null = $ps.PageHeight\n$null = $ps.TopMargin\n$null = $ps.BottomMargin\n$null = $ps.LeftMargin\n$null = $ps.RightMargin\n$null = $ps.HeaderDistance\n$null = $ps.FooterDistance\n$null = $ps.Gutter\n\nforeach ($s in $d.Styles) {\n    $null = $s.NameLocal\n}\n\n# No content, formatting, or style mutation is required: the template's\n# text, fields, bookmarks, section/page-setup values and style\n# definitions are unchanged by the source diff.\n"}
